$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "CNE" in A1 (previously A1 held a raw numeric id and had a
# bordered style inherited from the data column; the new header cell uses
# the default/no-border style like the other header cells).
$ws.Range("A1").Value = "CNE"
$ws.Range("A1").Borders.LineStyle = -4142

# B/C headers were swapped: B was "LastName"/C was "FirstName" before,
# now B is "FirstName"/C is "LastName".
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# The data in columns B and C (first/last name per student) were swapped
# to line up with the corrected headers.
for ($r = 2; $r -le 10; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    $cVal = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}

# Column A switched from auto-incrementing placeholder IDs (19000060...)
# to real CNE numbers (17000041...); only the first data row's value
# actually changes, the rest keep incrementing off of it via the existing
# formulas.
$ws.Range("A2").Value = 17000041

# A new (still-empty) row was appended below the table, inheriting the
# bordered style used by the rest of column A.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)

# Leave the selection where the author left off editing.
[void]$ws.Range("E9").Select()
